$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Cells.Item(17, 1).Value = "2025-07-24 15:42:11"
$ws.Cells.Item(17, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = "國立中山大學新海研3號貴重儀器使用中心誠徵專任技術員1名"
$ws.Cells.Item(17, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/ddc2e921-92c5-4004-8c2f-be2373c53f52?l=ch"
$ws.Cells.Item(17, 6).Value = "相關應徵資料予以保密，合者約談，不合者恕不另行通知。 發佈日期：2025-07-04 00:00:00"
$ws.Cells.Item(17, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(17, 8).Value = 300
$ws.Cells.Item(17, 9).Value = 30
$ws.Cells.Item(17, 10).Value = "hybrid_chunking"

# Row 18
$ws.Cells.Item(18, 1).Value = "2025-07-24 15:42:11"
$ws.Cells.Item(18, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = "[徵才] 國立臺灣大學防災減害與韌性學程 (綠‧韌性研究室) 徵求都市規劃/景觀/地理資訊專長 [專任計畫助理]"
$ws.Cells.Item(18, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/2793c7ef-b68d-4f00-9388-e011b78b9553?l=ch"
$ws.Cells.Item(18, 6).Value = "3.其他有利申請之相關文件 發佈日期：2025-07-21 00:00:00"
$ws.Cells.Item(18, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(18, 8).Value = 300
$ws.Cells.Item(18, 9).Value = 30
$ws.Cells.Item(18, 10).Value = "hybrid_chunking"

# Row 19
$ws.Cells.Item(19, 1).Value = "2025-07-24 15:42:11"
$ws.Cells.Item(19, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(19, 4).Value = "中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 博士後研究員"
$ws.Cells.Item(19, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/701ca4f1-a9f5-4a61-9b66-c4cf60f5c093?l=ch"
$ws.Cells.Item(19, 6).Value = "歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00"
$ws.Cells.Item(19, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(19, 8).Value = 300
$ws.Cells.Item(19, 9).Value = 30
$ws.Cells.Item(19, 10).Value = "hybrid_chunking"

# Row 20
$ws.Cells.Item(20, 1).Value = "2025-07-24 15:42:11"
$ws.Cells.Item(20, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = "中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 碩士級研究助理"
$ws.Cells.Item(20, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/2521ae27-55c0-4f27-9ded-b4bc908c1aff?l=ch"
$ws.Cells.Item(20, 6).Value = "歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00"
$ws.Cells.Item(20, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(20, 8).Value = 300
$ws.Cells.Item(20, 9).Value = 30
$ws.Cells.Item(20, 10).Value = "hybrid_chunking"

# Row 21
$ws.Cells.Item(21, 1).Value = "2025-07-24 15:42:11"
$ws.Cells.Item(21, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(21, 3).Value = 5
$ws.Cells.Item(21, 4).Value = "國立臺東大學通識教育中心徵聘專任助理教授以上教師徵才公告，收件至114年8月15日止。"
$ws.Cells.Item(21, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/e407fdbc-62c9-4e09-b08a-35a897cc4186?l=ch"
$ws.Cells.Item(21, 6).Value = "其    它： 相關訊息，請至本校首頁徵人啟事https://psn.nttu.edu.tw/p/406-1047-165359,r595.php?Lang=zh-tw查詢下載。 聯絡人姓名: 李家婕小姐 聯絡人電話: 089-517492 電子信箱：evalee@nttu.edu.tw 發佈日期：2025-07-09 00:00:00"
$ws.Cells.Item(21, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(21, 8).Value = 300
$ws.Cells.Item(21, 9).Value = 30
$ws.Cells.Item(21, 10).Value = "hybrid_chunking"

# Row 22
$ws.Cells.Item(22, 1).Value = "2025-07-24 16:04:57"
$ws.Cells.Item(22, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = "國立中山大學新海研3號貴重儀器使用中心誠徵專任技術員1名"
$ws.Cells.Item(22, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/ddc2e921-92c5-4004-8c2f-be2373c53f52?l=ch"
$ws.Cells.Item(22, 6).Value = "相關應徵資料予以保密，合者約談，不合者恕不另行通知。 發佈日期：2025-07-04 00:00:00"
$ws.Cells.Item(22, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(22, 8).Value = 300
$ws.Cells.Item(22, 9).Value = 30
$ws.Cells.Item(22, 10).Value = "hybrid_chunking"

# Row 23
$ws.Cells.Item(23, 1).Value = "2025-07-24 16:04:57"
$ws.Cells.Item(23, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = "[徵才] 國立臺灣大學防災減害與韌性學程 (綠‧韌性研究室) 徵求都市規劃/景觀/地理資訊專長 [專任計畫助理]"
$ws.Cells.Item(23, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/2793c7ef-b68d-4f00-9388-e011b78b9553?l=ch"
$ws.Cells.Item(23, 6).Value = "3.其他有利申請之相關文件 發佈日期：2025-07-21 00:00:00"
$ws.Cells.Item(23, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(23, 8).Value = 300
$ws.Cells.Item(23, 9).Value = 30
$ws.Cells.Item(23, 10).Value = "hybrid_chunking"

# Row 24
$ws.Cells.Item(24, 1).Value = "2025-07-24 16:04:57"
$ws.Cells.Item(24, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(24, 3).Value = 3
$ws.Cells.Item(24, 4).Value = "中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 博士後研究員"
$ws.Cells.Item(24, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/701ca4f1-a9f5-4a61-9b66-c4cf60f5c093?l=ch"
$ws.Cells.Item(24, 6).Value = "歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00"
$ws.Cells.Item(24, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(24, 8).Value = 300
$ws.Cells.Item(24, 9).Value = 30
$ws.Cells.Item(24, 10).Value = "hybrid_chunking"

# Row 25
$ws.Cells.Item(25, 1).Value = "2025-07-24 16:04:57"
$ws.Cells.Item(25, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(25, 3).Value = 4
$ws.Cells.Item(25, 4).Value = "中國醫藥大學 癌症生物精準醫學研究中心  王紹椿老師實驗室 誠徵 碩士級研究助理"
$ws.Cells.Item(25, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/2521ae27-55c0-4f27-9ded-b4bc908c1aff?l=ch"
$ws.Cells.Item(25, 6).Value = "歡迎對癌症研究有興趣的夥伴加入我們的團隊！ 發佈日期：2025-07-14 00:00:00"
$ws.Cells.Item(25, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(25, 8).Value = 300
$ws.Cells.Item(25, 9).Value = 30
$ws.Cells.Item(25, 10).Value = "hybrid_chunking"

# Row 26
$ws.Cells.Item(26, 1).Value = "2025-07-24 16:04:57"
$ws.Cells.Item(26, 2).Value = "材料相關的職缺有哪些？"
$ws.Cells.Item(26, 3).Value = 5
$ws.Cells.Item(26, 4).Value = "國立臺東大學通識教育中心徵聘專任助理教授以上教師徵才公告，收件至114年8月15日止。"
$ws.Cells.Item(26, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/e407fdbc-62c9-4e09-b08a-35a897cc4186?l=ch"
$ws.Cells.Item(26, 6).Value = "其    它： 相關訊息，請至本校首頁徵人啟事https://psn.nttu.edu.tw/p/406-1047-165359,r595.php?Lang=zh-tw查詢下載。 聯絡人姓名: 李家婕小姐 聯絡人電話: 089-517492 電子信箱：evalee@nttu.edu.tw 發佈日期：2025-07-09 00:00:00"
$ws.Cells.Item(26, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(26, 8).Value = 300
$ws.Cells.Item(26, 9).Value = 30
$ws.Cells.Item(26, 10).Value = "hybrid_chunking"

# Row 27
$ws.Cells.Item(27, 1).Value = "2025-07-24 16:05:29"
$ws.Cells.Item(27, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = "中國醫藥大學生物醫學研究所誠徵博士後研究員"
$ws.Cells.Item(27, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/cc8706e2-836d-4f85-94d2-89396360a823?l=ch"
$ws.Cells.Item(27, 6).Value = "生物醫學相關"
$ws.Cells.Item(27, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(27, 8).Value = 300
$ws.Cells.Item(27, 9).Value = 30
$ws.Cells.Item(27, 10).Value = "hybrid_chunking"

# Row 28
$ws.Cells.Item(28, 1).Value = "2025-07-24 16:05:29"
$ws.Cells.Item(28, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(28, 4).Value = "馬偕紀念醫院血液腫瘤科蘇迺文醫師誠徵國科會補助計畫專任助理"
$ws.Cells.Item(28, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/a78b7d93-b5b4-4bb9-bff5-5888b2d695e6?l=ch"
$ws.Cells.Item(28, 6).Value = "生命科學相關系所"
$ws.Cells.Item(28, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(28, 8).Value = 300
$ws.Cells.Item(28, 9).Value = 30
$ws.Cells.Item(28, 10).Value = "hybrid_chunking"

# Row 29
$ws.Cells.Item(29, 1).Value = "2025-07-24 16:05:29"
$ws.Cells.Item(29, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(29, 3).Value = 3
$ws.Cells.Item(29, 4).Value = "高雄榮總教研部生殖暨粒線體醫學研究室---誠徵博士後研究員"
$ws.Cells.Item(29, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/47729f59-955a-4b43-addd-5a18d1affa86?l=ch"
$ws.Cells.Item(29, 6).Value = "生物醫學相關領域畢業。"
$ws.Cells.Item(29, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(29, 8).Value = 300
$ws.Cells.Item(29, 9).Value = 30
$ws.Cells.Item(29, 10).Value = "hybrid_chunking"

# Row 30
$ws.Cells.Item(30, 1).Value = "2025-07-24 16:05:29"
$ws.Cells.Item(30, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(30, 3).Value = 4
$ws.Cells.Item(30, 4).Value = "台大醫院耳鼻喉部楊宗霖教授徵博士後研究員"
$ws.Cells.Item(30, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/e2f4f22d-3604-4ce0-854e-94b9a0ce8c10?l=ch"
$ws.Cells.Item(30, 6).Value = "細胞生物相關技術"
$ws.Cells.Item(30, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(30, 8).Value = 300
$ws.Cells.Item(30, 9).Value = 30
$ws.Cells.Item(30, 10).Value = "hybrid_chunking"

# Row 31
$ws.Cells.Item(31, 1).Value = "2025-07-24 16:05:29"
$ws.Cells.Item(31, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(31, 3).Value = 5
$ws.Cells.Item(31, 4).Value = "國家衛生研究院癌症研究所 誠徵院內博士後研究員或研究助理一名"
$ws.Cells.Item(31, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/95244417-fdb2-451f-812c-315ae9e234c4?l=ch"
$ws.Cells.Item(31, 6).Value = "具有生化、細胞、分生背景及細胞培養等相關研究經驗。"
$ws.Cells.Item(31, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(31, 8).Value = 300
$ws.Cells.Item(31, 9).Value = 30
$ws.Cells.Item(31, 10).Value = "hybrid_chunking"

# Row 32
$ws.Cells.Item(32, 1).Value = "2025-07-24 16:07:59"
$ws.Cells.Item(32, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(32, 4).Value = "中國醫藥大學生物醫學研究所誠徵博士後研究員"
$ws.Cells.Item(32, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/cc8706e2-836d-4f85-94d2-89396360a823?l=ch"
$ws.Cells.Item(32, 6).Value = "生物醫學相關"
$ws.Cells.Item(32, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(32, 8).Value = 300
$ws.Cells.Item(32, 9).Value = 30
$ws.Cells.Item(32, 10).Value = "hybrid_chunking"

# Row 33
$ws.Cells.Item(33, 1).Value = "2025-07-24 16:07:59"
$ws.Cells.Item(33, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(33, 3).Value = 2
$ws.Cells.Item(33, 4).Value = "馬偕紀念醫院血液腫瘤科蘇迺文醫師誠徵國科會補助計畫專任助理"
$ws.Cells.Item(33, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/a78b7d93-b5b4-4bb9-bff5-5888b2d695e6?l=ch"
$ws.Cells.Item(33, 6).Value = "生命科學相關系所"
$ws.Cells.Item(33, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(33, 8).Value = 300
$ws.Cells.Item(33, 9).Value = 30
$ws.Cells.Item(33, 10).Value = "hybrid_chunking"

# Row 34
$ws.Cells.Item(34, 1).Value = "2025-07-24 16:07:59"
$ws.Cells.Item(34, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(34, 3).Value = 3
$ws.Cells.Item(34, 4).Value = "高雄榮總教研部生殖暨粒線體醫學研究室---誠徵博士後研究員"
$ws.Cells.Item(34, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/47729f59-955a-4b43-addd-5a18d1affa86?l=ch"
$ws.Cells.Item(34, 6).Value = "生物醫學相關領域畢業。"
$ws.Cells.Item(34, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(34, 8).Value = 300
$ws.Cells.Item(34, 9).Value = 30
$ws.Cells.Item(34, 10).Value = "hybrid_chunking"

# Row 35
$ws.Cells.Item(35, 1).Value = "2025-07-24 16:07:59"
$ws.Cells.Item(35, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(35, 3).Value = 4
$ws.Cells.Item(35, 4).Value = "台大醫院耳鼻喉部楊宗霖教授徵博士後研究員"
$ws.Cells.Item(35, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/e2f4f22d-3604-4ce0-854e-94b9a0ce8c10?l=ch"
$ws.Cells.Item(35, 6).Value = "細胞生物相關技術"
$ws.Cells.Item(35, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(35, 8).Value = 300
$ws.Cells.Item(35, 9).Value = 30
$ws.Cells.Item(35, 10).Value = "hybrid_chunking"

# Row 36
$ws.Cells.Item(36, 1).Value = "2025-07-24 16:07:59"
$ws.Cells.Item(36, 2).Value = "生物相關的職缺有哪些"
$ws.Cells.Item(36, 3).Value = 5
$ws.Cells.Item(36, 4).Value = "國家衛生研究院癌症研究所 誠徵院內博士後研究員或研究助理一名"
$ws.Cells.Item(36, 5).Value = "https://www.nstc.gov.tw/folksonomy/detail/95244417-fdb2-451f-812c-315ae9e234c4?l=ch"
$ws.Cells.Item(36, 6).Value = "具有生化、細胞、分生背景及細胞培養等相關研究經驗。"
$ws.Cells.Item(36, 7).Value = "all-MiniLM-L6-v2"
$ws.Cells.Item(36, 8).Value = 300
$ws.Cells.Item(36, 9).Value = 30
$ws.Cells.Item(36, 10).Value = "hybrid_chunking"

Write-Output "Done writing rows 17-36"
